$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 3058
    7  = 1675
    10 = 35
    14 = 517
    15 = 350
    16 = 34
    23 = 3215
    25 = 139
    26 = 321
    29 = 95
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
